{"js": "// Update the division-problem worksheet table: replace the 25 problem\n// strings (5 populated rows x 5 columns; the interleaved rows are blank\n// answer rows and are left untouched) with their new values, in order,\n// preserving all existing run/paragraph formatting.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Row index within the table (0-based) -> new values for that row's 5 cells.\nconst newRowValues = {\n  0: [\"27\u00f76=\", \"21\u00f76=\", \"89\u00f79=\", \"47\u00f76=\", \"13\u00f77=\"],\n  4: [\"24\u00f76=\", \"91\u00f76=\", \"55\u00f78=\", \"99\u00f76=\", \"31\u00f78=\"],\n  8: [\"70\u00f74=\", \"53\u00f77=\", \"83\u00f76=\", \"88\u00f72=\", \"96\u00f72=\"],\n  12: [\"12\u00f77=\", \"17\u00f79=\", \"87\u00f78=\", \"33\u00f76=\", \"45\u00f72=\"],\n  16: [\"33\u00f79=\", \"69\u00f76=\", \"47\u00f78=\", \"75\u00f72=\", \"64\u00f76=\"],\n};\n\nfor (const rowIndexStr of Object.keys(newRowValues)) {\n  const rowIndex = Number(rowIndexStr);\n  const values = newRowValues[rowIndex];\n  for (let col = 0; col < values.length; col++) {\n    const cell = table.getCell(rowIndex, col);\n    cell.value = values[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the division-problem worksheet table: replace the 25 problem\n# strings (5 populated rows x 5 columns; the interleaved rows are blank\n# answer rows and are left untouched) with their new values, in order,\n# preserving all existing run/paragraph formatting.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @{\n    1  = @(\"27\u00f76=\", \"21\u00f76=\", \"89\u00f79=\", \"47\u00f76=\", \"13\u00f77=\")\n    5  = @(\"24\u00f76=\", \"91\u00f76=\", \"55\u00f78=\", \"99\u00f76=\", \"31\u00f78=\")\n    9  = @(\"70\u00f74=\", \"53\u00f77=\", \"83\u00f76=\", \"88\u00f72=\", \"96\u00f72=\")\n    13 = @(\"12\u00f77=\", \"17\u00f79=\", \"87\u00f78=\", \"33\u00f76=\", \"45\u00f72=\")\n    17 = @(\"33\u00f79=\", \"69\u00f76=\", \"47\u00f78=\", \"75\u00f72=\", \"64\u00f76=\")\n}\n\nforeach ($rowIndex in $newValues.Keys) {\n    $rowValues = $newValues[$rowIndex]\n    for ($col = 1; $col -le $rowValues.Length; $col++) {\n        $cell = $t.Cell($rowIndex, $col)\n        $cell.Range.Text = $rowValues[$col - 1]\n    }\n}\n"}
